$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 88 ("五感" entry) entirely; all rows below shift up by one.
$ws.Rows.Item(88).Delete()
